$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 0
$ws1.Range("F9").Value = 7146
$ws1.Range("F12").Value = 5495
$ws1.Range("F15").Value = 6271
$ws1.Range("F24").Value = 10153
$ws1.Range("F27").Value = 1976
$ws1.Range("F30").Value = 2132
$ws1.Range("F31").Value = 82
$ws1.Range("F32").Value = 0
$ws1.Range("F36").Value = 2088
$ws1.Range("F39").Value = 5228
$ws1.Range("F41").Value = 668
$ws1.Range("F46").Value = 986

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 9014
$ws4.Range("F11").Value = 7146
$ws4.Range("F14").Value = 5
$ws4.Range("F16").Value = 5495
$ws4.Range("F18").Value = 6271
$ws4.Range("F19").Value = 6271
$ws4.Range("F25").Value = 104
$ws4.Range("F27").Value = 10153
$ws4.Range("F30").Value = 1976
$ws4.Range("F32").Value = 2132
$ws4.Range("F33").Value = 82
$ws4.Range("F37").Value = 7
$ws4.Range("F38").Value = 2088
$ws4.Range("F40").Value = 5228
$ws4.Range("F42").Value = 668
$ws4.Range("F47").Value = 986
